$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-7 from 2023-09-01 (45170) to 2023-09-05 (45174)
$ws.Range("C2:C7").Value = 45174
